$d = $word.ActiveDocument

# The document currently has no styles part (word/styles.xml is absent).
# Materialize it by (re-)declaring the built-in "Normal" paragraph style,
# which is the only style referenced/used in this document.
$d.Styles.Add("Normal")
